# Apply updated cryptocurrency price/volume data to Sheet1 (D/E columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.385.01"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").Value = "1.979.92"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.79"
$ws.Range("E5").Value = "  -3.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  -3.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.43"
$ws.Range("E7").Value = "  -12.32%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.29"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +5.90%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.54"
$ws.Range("E13").Value = "  +8.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.859"
$ws.Range("E14").Value = "  -5.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.96"
$ws.Range("E15").Value = "  -6.90%  "
$ws.Range("D16").Value = "2.270.18"
$ws.Range("E16").Value = "  -3.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.43"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "1.979.87"
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").Value = "36.285.32"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.45"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("D21").Value = "0.0₃0860"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.31"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.06"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  -3.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.53"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.76"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  +7.85%  "
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.88"
$ws.Range("E33").Value = "  -7.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0630"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.40"
$ws.Range("E35").Value = "  -7.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.27"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.26"
$ws.Range("E38").Value = "  -7.43%  "
$ws.Range("E39").Value = "  -4.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.24"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0962"
$ws.Range("E42").Value = "  -6.69%  "
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0213"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.09"
$ws.Range("E45").Value = "  -5.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.21"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "92.21"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.55"
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("D49").Value = "1.370.88"
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.85"
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.96"
$ws.Range("E51").Value = "  -3.80%  "
